# epitweetr subscribers.xlsx template:
# 1. Making data directory read only
# 2. Add disclaimer on country / territories
# 3. Allow user to customize countries
#
# Spreadsheet-visible part of the change: insert a new "Real time Regions"
# header column right before the existing "Alert Slots" column, so the
# header row becomes:
#   User | Email | Topics | Excluded Topics | Real time Topics | Regions |
#   Real time Regions | Alert Slots

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Alert Slots" moves one column to the right, from H to I.
$ws.Cells.Item(2, 9).Value = $ws.Cells.Item(2, 8).Value()
$ws.Cells.Item(2, 9).Font.Bold = $true

# Column H now becomes the new "Real time Regions" header.
$ws.Cells.Item(2, 8).Value = "Real time Regions"
$ws.Cells.Item(2, 8).Font.Bold = $true

# Widen column H so the longer header text fits.
$ws.Columns.Item(8).ColumnWidth = 18.6

# Leave the selection on the newly edited header cell.
$ws.Range("H2").Select() | Out-Null
